# Reverse the order of the comma-separated "Recorded By" entries (column G)
# for every row whose value includes a "System" / "system" entry.
# Rows whose "Recorded By" value does not include "System" (e.g. two real
# user emails) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $value = $cell.Text

    if ($null -eq $value) { continue }
    if ($value -notlike "*,*") { continue }

    $parts = $value -split ",\s*"

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p -eq "System" -or $p -eq "system") {
            $hasSystem = $true
        }
    }

    if ($hasSystem) {
        $n = $parts.Length
        $reversed = @()
        for ($i = $n - 1; $i -ge 0; $i--) {
            $reversed += $parts[$i]
        }
        $cell.Value = $reversed -join ", "
    }
}
